# Add a new "2022-Q4" quarterly sheet right after the "总计" (summary) sheet,
# shifting the existing quarterly sheets down, and record the new quarter's
# totals on the summary sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Create the new sheet right after "总计" and name it "2022-Q4".
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add($null, $summary)
$newSheet.Name = "2022-Q4"

# NOTE: worksheet references must be (re-)fetched AFTER the sheet collection
# is mutated (Add/Delete/Move) - a reference obtained beforehand tracks the
# *positional slot*, not the sheet identity, and would now resolve to the
# freshly-inserted sheet instead of the one originally asked for.
$q3Sheet = $wb.Worksheets.Item("2022-Q3")

# Copy the overall cell formatting (fonts/borders/alignment) used throughout
# the other quarterly sheets so the new sheet matches the existing look:
#  - header row (B1:H1) uses the bold+bordered "header" style
#  - column A (the row-index column) uses the same bold+bordered style
$q3Sheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)   # xlPasteFormats
$q3Sheet.Range("A2").Copy()
$newSheet.Range("A2:A14").PasteSpecial(-4122)  # xlPasteFormats

# ---------------------------------------------------------------------------
# 2) Fill in the header row.
# ---------------------------------------------------------------------------
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# ---------------------------------------------------------------------------
# 3) Fill in the fund rows. Numeric-looking values in columns B, D, E, F, G
#    are written as TEXT (matching how this workbook stores formatted
#    numbers throughout), except where the value is exactly zero, which is
#    stored as a genuine number (same convention already used elsewhere in
#    the workbook).
# ---------------------------------------------------------------------------
$rows = @(
  @(0,  "007449", "兴全多维价值混合A",               "34.89", "86.54", "2.27", "0.7920",  8),
  @(1,  "007450", "兴全多维价值混合C",               "13.79", "86.54", "2.27", "0.3130",  8),
  @(2,  "014938", "同泰产业升级混合A",                "0.95", "68.98", "8.24", "0.0783",  1),
  @(3,  "004982", "新华安享多裕定期开放灵活配置混合",  "3.13", "45.81", "1.78", "0.0557", 10),
  @(4,  "014356", "长信企业成长三年持有混合A",         "1.70", "82.09", "2.83", "0.0481",  9),
  @(5,  "014357", "长信企业成长三年持有混合C",         "0.87", "82.09", "2.83", "0.0246",  9),
  @(6,  "008526", "华泰柏瑞行业精选混合A",             "0.56", "86.11", "2.49", "0.0139",  9),
  @(7,  "002409", "华夏新活力灵活配置混合A",           "0.12", "77.58", "8.46", "0.0102",  4),
  @(8,  "011361", "华夏博锐一年持有混合（MOM）A",      "0.11", "26.52", "4.49", "0.0049",  1),
  @(9,  "008527", "华泰柏瑞行业精选混合C",             "0.12", "86.11", "2.49", "0.0030",  9),
  @(10, "011362", "华夏博锐一年持有混合（MOM）C",      "0.00", "26.52", "4.49", $null,     1),
  @(11, "002410", "华夏新活力灵活配置混合C",           "0.00", "77.58", "8.46", $null,     4),
  @(12, "014939", "同泰产业升级混合C",                "-0.01", "68.98", "8.24", "-0.0008", 1)
)

$r = 2
foreach ($row in $rows) {
  $newSheet.Cells.Item($r, 1).Value = $row[0]

  $codeCell = $newSheet.Cells.Item($r, 2)
  $codeCell.NumberFormat = "@"
  $codeCell.Value = $row[1]

  $newSheet.Cells.Item($r, 3).Value = $row[2]

  $dCell = $newSheet.Cells.Item($r, 4)
  $dCell.NumberFormat = "@"
  $dCell.Value = $row[3]

  $eCell = $newSheet.Cells.Item($r, 5)
  $eCell.NumberFormat = "@"
  $eCell.Value = $row[4]

  $fCell = $newSheet.Cells.Item($r, 6)
  $fCell.NumberFormat = "@"
  $fCell.Value = $row[5]

  $gValue = $row[6]
  if ($null -eq $gValue) {
    $newSheet.Cells.Item($r, 7).Value = 0
  } else {
    $gCell = $newSheet.Cells.Item($r, 7)
    $gCell.NumberFormat = "@"
    $gCell.Value = $gValue
  }

  $newSheet.Cells.Item($r, 8).Value = $row[7]

  $r = $r + 1
}

# ---------------------------------------------------------------------------
# 4) Update the "总计" summary sheet: insert the 2022-Q4 totals as the new
#    first data row (row 2), pushing the other quarters down one row. The
#    index column (A) is just the 0-based row position, so every row is
#    simply rewritten with its final value.
# ---------------------------------------------------------------------------
$summary.Range("A5").Copy()
$summary.Range("A6").PasteSpecial(-4122)   # xlPasteFormats, extend style to new row

$summaryRows = @(
  @(0, "2022-Q4", 13, 1.34),
  @(1, "2022-Q3", 7,  0.25),
  @(2, "2022-Q2", 2,  0.02),
  @(3, "2021-Q4", 9,  0.93),
  @(4, "2021-Q3", 4,  0.39)
)

$r = 2
foreach ($row in $summaryRows) {
  $summary.Cells.Item($r, 1).Value = $row[0]
  $summary.Cells.Item($r, 2).Value = $row[1]
  $summary.Cells.Item($r, 3).Value = $row[2]
  $summary.Cells.Item($r, 4).Value = $row[3]
  $r = $r + 1
}
